# Updated cryptos list (price/volume refresh) matching the upstream data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.718.85"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "2.311.44"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.04"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.52"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.07"
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.53"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.40"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "2.661.81"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "2.314.72"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "42.686.14"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.33"
$ws.Range("E21").Value = "  +33.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.98"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.85"
$ws.Range("E24").Value = "  -4.93%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.68"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.88"
$ws.Range("E30").Value = "  +5.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.28"
$ws.Range("E31").Value = "  +7.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.74"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0890"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.58"
$ws.Range("E35").Value = "  -9.03%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("E40").Value = "  -5.83%  "
$ws.Range("E41").Value = "  +12.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.87"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.26"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.226"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.37"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.36"
$ws.Range("E47").Value = "  +3.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "81.97"
$ws.Range("E48").Value = "  +7.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.89"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "1.627.23"
$ws.Range("E51").Value = "  +4.97%  "
